# Weekly update: a new "Jengibre" (ginger) price record for
# "Vega Central Mapocho de Santiago" is inserted as the new row 19,
# pushing all the following rows (old 19..92) down by one (new 20..93).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 19 - this shifts rows 19:92 down to 20:93
# and (matching native Excel behaviour) carries the row-above's formatting
# onto the newly inserted row, which is exactly what the target file needs
# for the D column's date number-format (style index 2).
$ws.Rows("19").Insert()

# Populate the new row with this week's record.
$ws.Range("A19").Value = 9
$ws.Range("B19").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C19").Value = "Metropolitana"
$ws.Range("D19").Value = 44648
$ws.Range("E19").Value = 13
$ws.Range("F19").Value = 100114007
$ws.Range("G19").Value = "Jengibre"
$ws.Range("H19").Value = "Sin especificar"
$ws.Range("I19").Value = "Primera"
$ws.Range("J19").Value = 610
$ws.Range("K19").Value = 16000
$ws.Range("L19").Value = 17000
$ws.Range("M19").Value = 16500
$ws.Range("N19").Value = "$/caja 13 kilos"
$ws.Range("O19").Value = "Perú"
$ws.Range("P19").Value = 1269
$ws.Range("Q19").Value = 13
$ws.Range("R19").Value = "Hortaliza"
